$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A21 text changed from "compta" to "compta €" (style/format unchanged)
$ws.Range("A21").Value = "compta €"

# New rows 25-27: additional accounting format samples (euro FR, euro Irish, dollar Canada)
$ws.Range("A25").Value = "accounting € euro"
$ws.Range("B25").Value = 1.2
$ws.Range("B25").NumberFormat = '_-[$€-2]\ * #,##0.00_-;\-[$€-2]\ * #,##0.00_-;_-[$€-2]\ * "-"??_-;_-@_-'

$ws.Range("A26").Value = "accounting € Anglais Irlande"
$ws.Range("B26").Value = 2.3
$ws.Range("B26").NumberFormat = '_-[$€-1809]* #,##0.00_-;\-[$€-1809]* #,##0.00_-;_-[$€-1809]* "-"??_-;_-@_-'

$ws.Range("A27").Value = "accounting $ anglais - canada"
$ws.Range("B27").Value = 4.12
$ws.Range("B27").NumberFormat = '_-[$$-1009]* #,##0.00_-;\-[$$-1009]* #,##0.00_-;_-[$$-1009]* "-"??_-;_-@_-'

# Widen column A (to fit the new, longer labels) and move the selection to
# reflect the newly added last entry
$ws.Columns.Item(1).ColumnWidth = 30.6
$ws.Range("A27").Select()
